# "barchart buttons gefixt en initial page gemaakt"
# - Adds the "Amsterdam" data row (row 9) to the table.
# - Formats a set of helper/"button" cells (percentage style, right/center
#   aligned, Arial 8 - same look as the existing W22/X22/AG17:AK17 cells)
#   used to drive the bar-chart buttons further down the sheet.
# - Widens columns D and P so the new long labels fit.
# - Leaves the selection on D14, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New data row for "Amsterdam" (row 9)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Amsterdam"
$ws.Range("B9").Value = 0.32
$ws.Range("C9").Value = 0.46
$ws.Range("D9").Value = 0.22
$ws.Range("E9").Value = 0.41
$ws.Range("F9").Value = 0.33
$ws.Range("G9").Value = 0.11
$ws.Range("H9").Value = 0.14000000000000001
$ws.Range("I9").Value = 2489
$ws.Range("J9").Value = 0.32
$ws.Range("K9").Value = 0.16
$ws.Range("L9").Value = 0.14000000000000001
$ws.Range("M9").Value = 0.06
$ws.Range("N9").Value = 0.09
$ws.Range("O9").Value = 0.23
$ws.Range("P9").Value = 0.1
$ws.Range("Q9").Value = 0.16
$ws.Range("R9").Value = 0.22
$ws.Range("S9").Value = 0.35
$ws.Range("T9").Value = 0.13
$ws.Range("U9").Value = 0.05
$ws.Range("V9").Value = 0.06
$ws.Range("W9").Value = 0.18
$ws.Range("X9").Value = 0.11
$ws.Range("Y9").Value = 0.18
$ws.Range("Z9").Value = 0.28000000000000003
$ws.Range("AA9").Value = 0.19
$ws.Range("AB9").Value = 0.03
$ws.Range("AC9").Value = 0.28999999999999998
$ws.Range("AD9").Value = 0.33
$ws.Range("AE9").Value = 0.35

# ---------------------------------------------------------------------
# 2. "Button" cells - copy the existing button format (percentage,
#    Arial 8, right/center aligned - the same style already used by
#    W22/X22 and AG17:AK17) onto the newly-added button placeholders.
# ---------------------------------------------------------------------
$ws.Range("W22").Copy()

$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J11").PasteSpecial(-4122)
$ws.Range("J12").PasteSpecial(-4122)

$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("J13").PasteSpecial(-4122)
$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("Z13").PasteSpecial(-4122)

$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("Z14").PasteSpecial(-4122)

$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("Z15").PasteSpecial(-4122)

$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("Z16").PasteSpecial(-4122)

$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("O17").PasteSpecial(-4122)

$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("O18").PasteSpecial(-4122)

$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("S20").PasteSpecial(-4122)
$ws.Range("S21").PasteSpecial(-4122)
$ws.Range("S22").PasteSpecial(-4122)
$ws.Range("S23").PasteSpecial(-4122)
$ws.Range("S24").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Widen the label columns (D holds the long "button" captions, P the
#    income-bracket captions) so the new text isn't truncated.
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 23.5
$ws.Columns("P").ColumnWidth = 26

# ---------------------------------------------------------------------
# 4. Leave the selection where the author ended up while building the
#    initial page.
# ---------------------------------------------------------------------
[void]$ws.Range("D14").Select()
